$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the date strings in column I (shared strings) to ISO format (YYYY-MM-DD).
# NumberFormat stays General; a leading apostrophe keeps the literal text from being
# auto-parsed into a date serial, and resetting .Style afterwards drops the quote-prefix marker.
$c = $ws.Cells.Item(2, 9)
$c.Value = "'2020-03-19"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 9)
$c.Value = "'2020-03-19"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 9)
$c.Value = "'2020-03-20"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 9)
$c.Value = "'2020-03-20"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 9)
$c.Value = "'2020-03-21"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 9)
$c.Value = "'2020-03-21"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 9)
$c.Value = "'2020-03-22"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 9)
$c.Value = "'2020-03-22"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 9)
$c.Value = "'2020-03-23"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 9)
$c.Value = "'2020-03-23"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 9)
$c.Value = "'2020-03-24"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 9)
$c.Value = "'2020-03-24"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 9)
$c.Value = "'2020-03-25"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 9)
$c.Value = "'2020-03-25"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 9)
$c.Value = "'2020-03-26"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 9)
$c.Value = "'2020-03-26"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 9)
$c.Value = "'2020-03-27"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 9)
$c.Value = "'2020-03-27"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 9)
$c.Value = "'2020-03-28"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 9)
$c.Value = "'2020-03-28"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 9)
$c.Value = "'2020-03-29"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 9)
$c.Value = "'2020-03-29"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 9)
$c.Value = "'2020-03-30"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 9)
$c.Value = "'2020-03-30"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 9)
$c.Value = "'2020-03-31"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 9)
$c.Value = "'2020-03-31"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 9)
$c.Value = "'2020-04-01"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 9)
$c.Value = "'2020-04-01"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 9)
$c.Value = "'2020-04-02"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 9)
$c.Value = "'2020-04-02"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 9)
$c.Value = "'2020-04-03"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 9)
$c.Value = "'2020-04-03"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 9)
$c.Value = "'2020-04-04"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 9)
$c.Value = "'2020-04-04"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 9)
$c.Value = "'2020-04-05"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 9)
$c.Value = "'2020-04-05"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 9)
$c.Value = "'2020-04-06"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 9)
$c.Value = "'2020-04-06"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 9)
$c.Value = "'2020-04-07"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 9)
$c.Value = "'2020-04-07"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 9)
$c.Value = "'2020-04-08"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 9)
$c.Value = "'2020-04-08"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 9)
$c.Value = "'2020-04-09"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 9)
$c.Value = "'2020-04-09"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 9)
$c.Value = "'2020-04-10"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 9)
$c.Value = "'2020-04-10"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 9)
$c.Value = "'2020-04-11"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 9)
$c.Value = "'2020-04-11"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 9)
$c.Value = "'2020-04-12"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 9)
$c.Value = "'2020-04-12"
$c.Style = "Normal"
$c = $ws.Cells.Item(52, 9)
$c.Value = "'2020-04-13"
$c.Style = "Normal"
$c = $ws.Cells.Item(53, 9)
$c.Value = "'2020-04-13"
$c.Style = "Normal"
$c = $ws.Cells.Item(54, 9)
$c.Value = "'2020-04-14"
$c.Style = "Normal"
$c = $ws.Cells.Item(55, 9)
$c.Value = "'2020-04-14"
$c.Style = "Normal"
$c = $ws.Cells.Item(56, 9)
$c.Value = "'2020-04-15"
$c.Style = "Normal"
$c = $ws.Cells.Item(57, 9)
$c.Value = "'2020-04-15"
$c.Style = "Normal"
$c = $ws.Cells.Item(58, 9)
$c.Value = "'2020-04-16"
$c.Style = "Normal"
$c = $ws.Cells.Item(59, 9)
$c.Value = "'2020-04-16"
$c.Style = "Normal"
$c = $ws.Cells.Item(60, 9)
$c.Value = "'2020-04-17"
$c.Style = "Normal"
$c = $ws.Cells.Item(61, 9)
$c.Value = "'2020-04-17"
$c.Style = "Normal"
$c = $ws.Cells.Item(62, 9)
$c.Value = "'2020-04-18"
$c.Style = "Normal"
$c = $ws.Cells.Item(63, 9)
$c.Value = "'2020-04-18"
$c.Style = "Normal"
$c = $ws.Cells.Item(64, 9)
$c.Value = "'2020-04-19"
$c.Style = "Normal"
$c = $ws.Cells.Item(65, 9)
$c.Value = "'2020-04-19"
$c.Style = "Normal"
$c = $ws.Cells.Item(66, 9)
$c.Value = "'2020-04-20"
$c.Style = "Normal"
$c = $ws.Cells.Item(67, 9)
$c.Value = "'2020-04-20"
$c.Style = "Normal"
$c = $ws.Cells.Item(68, 9)
$c.Value = "'2020-04-21"
$c.Style = "Normal"
$c = $ws.Cells.Item(69, 9)
$c.Value = "'2020-04-21"
$c.Style = "Normal"
$c = $ws.Cells.Item(70, 9)
$c.Value = "'2020-04-22"
$c.Style = "Normal"
$c = $ws.Cells.Item(71, 9)
$c.Value = "'2020-04-22"
$c.Style = "Normal"
$c = $ws.Cells.Item(72, 9)
$c.Value = "'2020-04-23"
$c.Style = "Normal"
$c = $ws.Cells.Item(73, 9)
$c.Value = "'2020-04-23"
$c.Style = "Normal"
$c = $ws.Cells.Item(74, 9)
$c.Value = "'2020-04-24"
$c.Style = "Normal"
$c = $ws.Cells.Item(75, 9)
$c.Value = "'2020-04-24"
$c.Style = "Normal"
$c = $ws.Cells.Item(76, 9)
$c.Value = "'2020-04-25"
$c.Style = "Normal"
$c = $ws.Cells.Item(77, 9)
$c.Value = "'2020-04-25"
$c.Style = "Normal"
$c = $ws.Cells.Item(78, 9)
$c.Value = "'2020-04-26"
$c.Style = "Normal"
$c = $ws.Cells.Item(79, 9)
$c.Value = "'2020-04-26"
$c.Style = "Normal"
$c = $ws.Cells.Item(80, 9)
$c.Value = "'2020-04-27"
$c.Style = "Normal"
$c = $ws.Cells.Item(81, 9)
$c.Value = "'2020-04-27"
$c.Style = "Normal"
$c = $ws.Cells.Item(82, 9)
$c.Value = "'2020-04-28"
$c.Style = "Normal"
$c = $ws.Cells.Item(83, 9)
$c.Value = "'2020-04-28"
$c.Style = "Normal"
$c = $ws.Cells.Item(84, 9)
$c.Value = "'2020-04-29"
$c.Style = "Normal"
$c = $ws.Cells.Item(85, 9)
$c.Value = "'2020-04-29"
$c.Style = "Normal"
$c = $ws.Cells.Item(86, 9)
$c.Value = "'2020-04-30"
$c.Style = "Normal"
$c = $ws.Cells.Item(87, 9)
$c.Value = "'2020-04-30"
$c.Style = "Normal"
$c = $ws.Cells.Item(88, 9)
$c.Value = "'2020-05-01"
$c.Style = "Normal"
$c = $ws.Cells.Item(89, 9)
$c.Value = "'2020-05-01"
$c.Style = "Normal"
$c = $ws.Cells.Item(90, 9)
$c.Value = "'2020-05-02"
$c.Style = "Normal"
$c = $ws.Cells.Item(91, 9)
$c.Value = "'2020-05-02"
$c.Style = "Normal"
$c = $ws.Cells.Item(92, 9)
$c.Value = "'2020-05-03"
$c.Style = "Normal"
$c = $ws.Cells.Item(93, 9)
$c.Value = "'2020-05-03"
$c.Style = "Normal"
$c = $ws.Cells.Item(94, 9)
$c.Value = "'2020-05-04"
$c.Style = "Normal"
$c = $ws.Cells.Item(95, 9)
$c.Value = "'2020-05-04"
$c.Style = "Normal"

# Column A (the old pandas reset_index column) becomes a plain running count.
for ($r = 4; $r -le 95; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Append the two new rows for 2020-05-05 (Kewaunee + Manitowoc), matching the existing pattern.
$ws.Cells.Item(95, 1).Copy()
$ws.Cells.Item(96, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(96, 1).Value = 94
$ws.Cells.Item(96, 2).Value = 55061
$ws.Cells.Item(96, 3).Value = "Point Beach"
$ws.Cells.Item(96, 4).Value = 3
$ws.Cells.Item(96, 5).Value = "Kewaunee"
$ws.Cells.Item(96, 6).Value = "Wisconsin"
$ws.Cells.Item(96, 7).Value = "Kewaunee, Wisconsin, US"
$ws.Cells.Item(96, 8).Value = 22
$c = $ws.Cells.Item(96, 9)
$c.Value = "'2020-05-05"
$c.Style = "Normal"

$ws.Cells.Item(95, 1).Copy()
$ws.Cells.Item(97, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(97, 1).Value = 95
$ws.Cells.Item(97, 2).Value = 55071
$ws.Cells.Item(97, 3).Value = "Point Beach"
$ws.Cells.Item(97, 4).Value = 3
$ws.Cells.Item(97, 5).Value = "Manitowoc"
$ws.Cells.Item(97, 6).Value = "Wisconsin"
$ws.Cells.Item(97, 7).Value = "Manitowoc, Wisconsin, US"
$ws.Cells.Item(97, 8).Value = 16
$c = $ws.Cells.Item(97, 9)
$c.Value = "'2020-05-05"
$c.Style = "Normal"

